$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = new price text; G = new "Hora" text }
# Values are entered with a leading apostrophe so Excel stores them as
# literal text (preserving formats like trailing zeros, e.g. "0.06197"),
# matching the original inline-string cell contents instead of coercing
# them into numbers.
$updates = @{
    2 = @{ D="272.75"; G="2" }
    3 = @{ D="21.10"; G="2" }
    4 = @{ D="6.264"; G="2" }
    5 = @{ D="0.06197"; G="2" }
    6 = @{ D="3.560"; G="2" }
    7 = @{ D="6.544"; G="2" }
    8 = @{ D="1.468"; G="2" }
    9 = @{ D="0.8268"; G="2" }
    10 = @{ D="0.1663"; G="2" }
    11 = @{ D="0.08274"; G="2" }
    12 = @{ D="0.03522"; G="2" }
    13 = @{ D="0.03191"; G="2" }
    14 = @{ D="0.09181"; G="2" }
    15 = @{ D="3.766"; G="2" }
    16 = @{ D="0.001648"; G="2" }
    17 = @{ D="0.04712"; G="2" }
    18 = @{ D="0.006385"; G="2" }
    19 = @{ D="0.006196"; G="2" }
    20 = @{ G="2" }
    21 = @{ D="0.0001501"; G="2" }
    22 = @{ D="3.733"; G="2" }
    23 = @{ D="2.257"; G="2" }
    24 = @{ D="0.01386"; G="2" }
    25 = @{ D="0.3292"; G="2" }
    26 = @{ G="2" }
    27 = @{ G="2" }
    28 = @{ D="0.0002715"; G="2" }
    29 = @{ G="2" }
    30 = @{ G="2" }
    31 = @{ G="2" }
    32 = @{ G="2" }
    33 = @{ G="2" }
    34 = @{ G="2" }
    35 = @{ G="2" }
    36 = @{ G="2" }
    37 = @{ G="2" }
    38 = @{ G="2" }
    39 = @{ G="2" }
    40 = @{ D="0.04701"; G="2" }
    41 = @{ D="0.007017"; G="2" }
    42 = @{ D="0.004102"; G="2" }
    43 = @{ D="0.1117"; G="2" }
    44 = @{ D="0.01038"; G="2" }
    45 = @{ D="0.00006304"; G="2" }
    46 = @{ D="0.0009905"; G="2" }
    47 = @{ G="2" }
    48 = @{ D="0.9506"; G="2" }
    49 = @{ D="0.001399"; G="2" }
    50 = @{ D="0.00001901"; G="2" }
    51 = @{ D="0.01241"; G="2" }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($u.ContainsKey("D")) {
        $ws.Range("D$row").Value = "'" + $u.D
    }
    $ws.Range("G$row").Value = "'" + $u.G
}
